$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The report lists low-stock items alphabetically (rows 7-20). Two new
# items need to be inserted in their correct alphabetical position
# (between "DAVALINDI..." at row 8 and "DOLIPRANE..." at row 9):
#   - DIASTOP SUSP. 60ML
#   - DIVIDO 75MG 30 DUAL RELEASE CAPS.
# Inserting 2 full rows at row 9 pushes everything below (including the
# totals row and the footer row) down by two rows, same as a normal
# Excel "Insert Rows" operation.
$ws.Rows("9:10").Insert()

# Restore the row heights for the two new rows (Insert() left them at
# the default/blank height).
$ws.Rows("9:9").RowHeight = 25.5
$ws.Rows("10:10").RowHeight = 24.75

# Recreate the same merge pattern used by every other item row.
$ws.Range("A9:B9").Merge()
$ws.Range("C9:G9").Merge()
$ws.Range("H9:K9").Merge()
$ws.Range("L9:M9").Merge()
$ws.Range("N9:O9").Merge()

$ws.Range("A10:B10").Merge()
$ws.Range("C10:G10").Merge()
$ws.Range("H10:K10").Merge()
$ws.Range("L10:M10").Merge()
$ws.Range("N10:O10").Merge()

# Row 9: DIASTOP SUSP. 60ML
# (L and P are formatted with numeric formats even though the report
# stores them as plain text, so a leading apostrophe is used to force
# text entry instead of letting Excel auto-convert them to numbers.)
$ws.Range("A9").Value = 3
$ws.Range("C9").Value = "DIASTOP SUSP. 60ML"
$ws.Range("H9").Value = "0:0"
$ws.Range("L9").Value = "'1"
$ws.Range("N9").Value = "30.00"
$ws.Range("P9").Value = "'30.0000"
$ws.Range("Q9").Value = "1:0"

# Row 10: DIVIDO 75MG 30 DUAL RELEASE CAPS.
$ws.Range("A10").Value = 4
$ws.Range("C10").Value = "DIVIDO 75MG 30 DUAL RELEASE CAPS."
$ws.Range("H10").Value = "3:0"
$ws.Range("L10").Value = "'1"
$ws.Range("N10").Value = "141.00"
$ws.Range("P10").Value = "'46.5300"
$ws.Range("Q10").Value = "0:1"

# Renumber the "م" (item #) column for all the rows that followed the
# insertion point, since they've shifted down by two rows but still
# need to count up sequentially from 1.
$ws.Range("A11").Value = 5
$ws.Range("A12").Value = 6
$ws.Range("A13").Value = 7
$ws.Range("A14").Value = 8
$ws.Range("A15").Value = 9
$ws.Range("A16").Value = 10
$ws.Range("A17").Value = 11
$ws.Range("A18").Value = 12
$ws.Range("A19").Value = 13
$ws.Range("A20").Value = 14
$ws.Range("A21").Value = 15
$ws.Range("A22").Value = 16

# Update the total (sum of the "sale price" column) to include the two
# newly-added rows.
$ws.Range("P23").Value = 811.17

# The report was regenerated a few minutes later, so the timestamp in
# the footer advances from 11:07 AM to 11:16 AM.
$ws.Range("A24").Value = "Sunday, 27 July, 2025 11:16 AM"
